# Mongo lab creds updated: add a "Participant Name" column before the
# existing Username/Password columns and populate it with the roster.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new column at the front; everything currently in A:G
# (Username, Password, ...) shifts right into B:H.
$ws.Columns.Item(1).EntireColumn.Insert()

$names = @(
    "Ashok Kumar",
    "Bhagyasree.N",
    "Bhargavi Bhargavi",
    "Chaandan Banerjee",
    "Deepak Naita",
    "Eknath Vashishtha",
    "Hemalatha Enugu",
    "Hritvik Dekate",
    "Karthika Rajaram",
    "Kumari Divya",
    "Manoj Pradhan",
    "Megha R",
    "Modugula Supriya",
    "Nandini S",
    "Nikhil Dhaka",
    "Pankaj Rawat",
    "Piyush Goyal",
    "Sandhya Shiramagond",
    "Sathish Kumar K",
    "SATHISH RANGAN",
    "Sivaraman L",
    "Sriram Sarveswaran",
    "Subham",
    "Subhani Shaik",
    "Sunil Pallath Sagar"
)

$ws.Range("A1").Value = "Participant Name"
for ($i = 0; $i -lt $names.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
}

Write-Host "Done"
